# Combine Florenceville and Bristol prior to 2009
#
# Starting with 2009, Florenceville and Bristol were amalgamated into
# Florenceville-Bristol. For data consistency, the two standalone rows
# ("Bristol" and "Florenceville") are removed from this pre-2009 sheet,
# and the Policing Provider values for the remaining municipalities are
# corrected to reflect the manually-combined/weighted source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the two rows that no longer exist on their own ---------

# "Bristol" is row 16.
if ($ws.Range("A16").Text -ne "Bristol") {
    throw "Expected row 16 to be Bristol, found '$($ws.Range('A16').Text)'"
}
$ws.Rows(16).Delete()

# "Florenceville" was row 33, now row 32 after the previous deletion.
if ($ws.Range("A32").Text -ne "Florenceville") {
    throw "Expected row 32 to be Florenceville, found '$($ws.Range('A32').Text)'"
}
$ws.Rows(32).Delete()

# --- 2. Correct the Policing Provider column for the affected rows ----
# (row numbers below are final, i.e. after the two deletions above)

function Set-Provider($row, $municipality, $provider) {
    $actual = $ws.Range("A$row").Text
    if ($actual -ne $municipality) {
        throw "Row $row : expected '$municipality', found '$actual'"
    }
    $ws.Range("B$row").Value = $provider
}

Set-Provider 16  "Cambridge-Narrows"           "MPSA"
Set-Provider 17  "Campbellton"                 "PPSA"
Set-Provider 25  "Dalhousie"                   "MPSA"
Set-Provider 26  "Dieppe"                      "PPSA"
Set-Provider 27  "Doaktown"                    "MPSA"
Set-Provider 28  "Dorchester"                  "Municipal"
Set-Provider 30  "Edmundston"                  "PPSA"
Set-Provider 33  "Fredericton"                 "PPSA"
Set-Provider 35  "Gagetown"                    "Municipal"
Set-Provider 37  "Grand Falls/Grand-Sault"     "PPSA"
Set-Provider 38  "Grand Manan"                 "MPSA"
Set-Provider 40  "Hampton"                     "PPSA"
Set-Provider 52  "Millville"                   "Municipal"
Set-Provider 53  "Minto"                       "MPSA"
Set-Provider 54  "Miramichi"                   "PPSA"
Set-Provider 55  "Moncton"                     "PPSA"
Set-Provider 56  "Nackawic"                    "Municipal"
Set-Provider 58  "Nigadoo"                     "PPSA"
Set-Provider 59  "Norton"                      "MPSA"
Set-Provider 61  "Oromocto"                    "PPSA"
Set-Provider 62  "Paquetville"                 "Municipal"
Set-Provider 64  "Petit-Rocher"                "PPSA"
Set-Provider 65  "Petitcodiac"                 "Municipal"
Set-Provider 69  "Quispamsis"                  "PPSA"
Set-Provider 71  "Richibucto"                  "MPSA"
Set-Provider 72  "Riverside-Albert"            "Municipal"
Set-Provider 73  "Riverview"                   "PPSA"
Set-Provider 75  "Rogersville"                 "MPSA"
Set-Provider 76  "Rothesay"                    "MPSA"
Set-Provider 77  "Sackville"                   "PPSA"
Set-Provider 78  "Saint Andrews"               "Municipal"
Set-Provider 79  "Saint George"                "Municipal"
Set-Provider 80  "Saint John"                  "PPSA"
Set-Provider 81  "Saint-André"                 "PPSA"
Set-Provider 98  "Sussex"                      "MPSA"
Set-Provider 100 "Tide Head"                   "PPSA"
